# Apply the "2021-05 Victorian Outbreak Paths" update:
#  - Sheet1: fix the F68/F71 link-label text, insert one new contact-tracing
#    row (West Melbourne / Student) before the old row 78, and append one new
#    row (Arcare Maidstone / Household) at the end of Table1.
#  - "Date Colours" sheet: refresh the 14-colour palette in column B and add
#    a new Date/Colour row to Date_Colours.
#  - Leave Sheet1 as the active sheet/selection afterwards.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Date Colours")

# ---------------------------------------------------------------------------
# 1. Sheet1: correct the existing "Student, North Melbourne Primary School"
#    link label down to plain "Student" (rows 68 & 71).
# ---------------------------------------------------------------------------
$ws1.Range("F68").Value = "Student"
$ws1.Range("F71").Value = "Student"

# ---------------------------------------------------------------------------
# 2. Sheet1 / Table1: insert a new row above the current row 78 (shifts the
#    existing rows 78-86 down to 79-87), then populate it.
# ---------------------------------------------------------------------------
$ws1.Rows("78:78").Insert()

$ws1.Range("A78").Value = 44355
$ws1.Range("A78").NumberFormat = $ws1.Range("A77").NumberFormat
$ws1.Range("B78").Value = "M # b"
$ws1.Range("C78").Value = "M # n"
$ws1.Range("D78").Value = "West Melbourne"
$ws1.Range("F78").Value = "Student"
$ws1.Range("G78").Value = "Delta (B.1.617.2)"

# ---------------------------------------------------------------------------
# 3. Sheet1 / Table1: append a new row (row 88) after the old last row
#    (now row 87).
# ---------------------------------------------------------------------------
$ws1.Range("A88").Value = 44355
$ws1.Range("A88").NumberFormat = $ws1.Range("A87").NumberFormat
$ws1.Range("B88").Value = "A # a"
$ws1.Range("C88").Value = "A # d"
$ws1.Range("D88").Value = "Arcare Maidstone"
$ws1.Range("F88").Value = "Household"
$ws1.Range("G88").Value = "Kappa (B.1.617.1)"

$tbl1 = $ws1.ListObjects.Item(1)
$tbl1.Resize($ws1.Range("A1:G88"))

# ---------------------------------------------------------------------------
# 4. "Date Colours" sheet / Date_Colours table: refresh the colour palette
#    used in column B (rows 2-15).
# ---------------------------------------------------------------------------
$ws2.Range("B2").Value = "#f5f4ff"
$ws2.Range("B3").Value = "#ece9fe"
$ws2.Range("B4").Value = "#e2defd"
$ws2.Range("B5").Value = "#d8d4fd"
$ws2.Range("B6").Value = "#cdc9fc"
$ws2.Range("B7").Value = "#c3bffb"
$ws2.Range("B8").Value = "#b8b4fa"
$ws2.Range("B9").Value = "#adaaf9"
$ws2.Range("B10").Value = "#a2a0f8"
$ws2.Range("B11").Value = "#9696f7"
$ws2.Range("B12").Value = "#898cf6"
$ws2.Range("B13").Value = "#7c82f5"
$ws2.Range("B14").Value = "#6d79f4"
$ws2.Range("B15").Value = "#5d6ff2"

# ---------------------------------------------------------------------------
# 5. "Date Colours" / Date_Colours: append a new Date/Colour row (row 16).
# ---------------------------------------------------------------------------
$ws2.Range("A16").Value = 44355
$ws2.Range("A16").NumberFormat = $ws2.Range("A15").NumberFormat
$ws2.Range("B16").Value = "#4966f1"
$ws2.Range("C16").Value = "#CC66FF"

$tbl2 = $ws2.ListObjects.Item(1)
$tbl2.Resize($ws2.Range("A1:C16"))

# ---------------------------------------------------------------------------
# 6. Make Sheet1 the active sheet/tab again, with A88 selected (the newly
#    appended row), matching the saved view state.
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A88").Select()
